# Auto-generated edit script applying numeric value updates
# as described by the commit diff (refreshed market-price data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2572.6
$ws.Range("I19").Value = 4350
$ws.Range("K19").Value = 4350
$ws.Range("M19").Value = -4175

$ws.Range("H43").Value = 891.7692
$ws.Range("I43").Value = 549.75
$ws.Range("J43").Value = 1043.7778
$ws.Range("K43").Value = 549.75
$ws.Range("L43").Value = 1043.7778
$ws.Range("M43").Value = -480.75
$ws.Range("N43").Value = -1181.7778

$ws.Range("H113").Value = 41670200
$ws.Range("I113").Value = 55558196
$ws.Range("J113").Value = 6216.1665
$ws.Range("K113").Value = 55558196
$ws.Range("L113").Value = 6216.1665
$ws.Range("M113").Value = -55554942
$ws.Range("N113").Value = -12724.1665

$ws.Range("H116").Value = 3907
$ws.Range("I116").Value = 1803.3636
$ws.Range("K116").Value = 1803.3636
$ws.Range("M116").Value = 1638.6364

$ws.Range("H129").Value = 223167.33
$ws.Range("J129").Value = 278904.22
$ws.Range("L129").Value = 836712.6599999999
$ws.Range("N129").Value = -846712.6599999999

$ws.Range("H132").Value = 3201.4849
$ws.Range("I132").Value = 3477.4827
$ws.Range("J132").Value = 1200.5
$ws.Range("K132").Value = 10432.4481
$ws.Range("L132").Value = 3601.5
$ws.Range("M132").Value = -7902.4481
$ws.Range("N132").Value = -8661.5

$ws.Range("H138").Value = 30305144
$ws.Range("J138").Value = 3183.4285
$ws.Range("L138").Value = 9550.2855
$ws.Range("N138").Value = -19830.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 708.12
$ws.Range("I2").Value = 788.3889
$ws.Range("J2").Value = 501.7143
$ws.Range("K2").Value = 788.3889
$ws.Range("L2").Value = 501.7143
$ws.Range("M2").Value = -675.3889
$ws.Range("N2").Value = -727.7143

$ws.Range("H32").Value = 6960.6807
$ws.Range("I32").Value = 4907.755
$ws.Range("K32").Value = 4907.755
$ws.Range("M32").Value = -4620.755

$ws.Range("H45").Value = 2560.0952
$ws.Range("I45").Value = 1776.7693
$ws.Range("J45").Value = 3833
$ws.Range("K45").Value = 1776.7693
$ws.Range("L45").Value = 3833
$ws.Range("M45").Value = -1399.7693
$ws.Range("N45").Value = -4587

$ws.Range("H74").Value = 28572820
$ws.Range("I74").Value = 38461964
$ws.Range("K74").Value = 38461964
$ws.Range("M74").Value = -38461090

$ws.Range("H77").Value = 28572820
$ws.Range("I77").Value = 38461964
$ws.Range("K77").Value = 192309820
$ws.Range("M77").Value = -192305452

$ws.Range("H102").Value = 1668.7142
$ws.Range("I102").Value = 1617.5
$ws.Range("J102").Value = 1737
$ws.Range("K102").Value = 1617.5
$ws.Range("L102").Value = 1737
$ws.Range("M102").Value = 4.5
$ws.Range("N102").Value = -4981

$ws.Range("H110").Value = 758.3
$ws.Range("I110").Value = 647.875
$ws.Range("K110").Value = 647.875
$ws.Range("M110").Value = 1397.125

$ws.Range("H116").Value = 708.12
$ws.Range("I116").Value = 788.3889
$ws.Range("J116").Value = 501.7143
$ws.Range("K116").Value = 788.3889
$ws.Range("L116").Value = 501.7143
$ws.Range("M116").Value = 1505.6111
$ws.Range("N116").Value = -5089.7143

$ws.Range("H122").Value = 3604.8572
$ws.Range("I122").Value = 2539.1667
$ws.Range("K122").Value = 7617.500100000001
$ws.Range("M122").Value = -5167.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 708.12
$ws.Range("I3").Value = 788.3889
$ws.Range("J3").Value = 501.7143
$ws.Range("K3").Value = 788.3889
$ws.Range("L3").Value = 501.7143
$ws.Range("M3").Value = -674.3889
$ws.Range("N3").Value = -729.7143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 18672.379
$ws.Range("I58").Value = 1294.6316
$ws.Range("J58").Value = 51690.1
$ws.Range("K58").Value = 1294.6316
$ws.Range("L58").Value = 51690.1
$ws.Range("M58").Value = -1091.6316
$ws.Range("N58").Value = -52096.1

$ws.Range("H107").Value = 1084.5385
$ws.Range("I107").Value = 381.4375
$ws.Range("J107").Value = 2209.5
$ws.Range("K107").Value = 381.4375
$ws.Range("L107").Value = 2209.5
$ws.Range("M107").Value = 1538.5625
$ws.Range("N107").Value = -6049.5

$ws.Range("H132").Value = 4728.5713
$ws.Range("I132").Value = 3703
$ws.Range("K132").Value = 11109
$ws.Range("M132").Value = -8579

$ws.Range("H134").Value = 1857.1428
$ws.Range("I134").Value = 1800
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5400
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -2865
$ws.Range("N134").Value = -11070

$ws.Range("H136").Value = 18672.379
$ws.Range("I136").Value = 1294.6316
$ws.Range("J136").Value = 51690.1
$ws.Range("K136").Value = 3883.8948
$ws.Range("L136").Value = 155070.3
$ws.Range("M136").Value = -1333.8948
$ws.Range("N136").Value = -160170.3

$ws.Range("H137").Value = 30640
$ws.Range("J137").Value = 30640
$ws.Range("L137").Value = 30640
$ws.Range("N137").Value = -40840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 294.36365
$ws.Range("I14").Value = 294.36365
$ws.Range("K14").Value = 883.09095
$ws.Range("M14").Value = -710.09095

$ws.Range("H17").Value = 401
$ws.Range("I17").Value = 314.75
$ws.Range("K17").Value = 944.25
$ws.Range("M17").Value = -775.25

$ws.Range("H131").Value = 737.59
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 742.0101
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2226.0303
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12306.0303

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H46").Value = 16000
$ws.Range("J46").Value = 12000
$ws.Range("L46").Value = 12000
$ws.Range("N46").Value = -12312

$ws.Range("H107").Value = 3344749.2
$ws.Range("I107").Value = 279.7857
$ws.Range("K107").Value = 279.7857
$ws.Range("M107").Value = 1640.2143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1006.8485
$ws.Range("I46").Value = 990.86664
$ws.Range("J46").Value = 1166.6666
$ws.Range("K46").Value = 990.86664
$ws.Range("L46").Value = 1166.6666
$ws.Range("M46").Value = -802.86664
$ws.Range("N46").Value = -1542.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 19500
$ws.Range("J69").Value = 19500
$ws.Range("L69").Value = 19500
$ws.Range("N69").Value = -20998

$ws.Range("H72").Value = 19500
$ws.Range("J72").Value = 19500
$ws.Range("L72").Value = 58500
$ws.Range("N72").Value = -65988

$ws.Range("H81").Value = 166668660
$ws.Range("I81").Value = 1675
$ws.Range("J81").Value = 500002620
$ws.Range("K81").Value = 3350
$ws.Range("L81").Value = 1000005240
$ws.Range("M81").Value = -2289
$ws.Range("N81").Value = -1000007362

$ws.Range("H84").Value = 166668660
$ws.Range("I84").Value = 1675
$ws.Range("J84").Value = 500002620
$ws.Range("K84").Value = 16750
$ws.Range("L84").Value = 5000026200
$ws.Range("M84").Value = -11446
$ws.Range("N84").Value = -5000036808

$ws.Range("H122").Value = 1155.7333
$ws.Range("I122").Value = 1044.0454
$ws.Range("J122").Value = 1462.875
$ws.Range("K122").Value = 3132.1362
$ws.Range("L122").Value = 4388.625
$ws.Range("M122").Value = -682.1361999999999
$ws.Range("N122").Value = -9288.625

$ws.Range("H136").Value = 25808014
$ws.Range("I136").Value = 32259288
$ws.Range("K136").Value = 96777864
$ws.Range("M136").Value = -96775314
